$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 8 de Octubre de 2020 a las 12:45
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 12:45"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7776796
$ws.Range("C4").Value = 572
$ws.Range("D4").Value = 4984154
$ws.Range("E4").Value = 2575854
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 216788

# Row 19: Banglades
$ws.Range("A19").Value = "Banglades"
$ws.Range("B19").Value = 374592
$ws.Range("C19").Value = 1441
$ws.Range("D19").Value = 288316
$ws.Range("E19").Value = 80816
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 5460

# Row 31: Rumania
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 145700
$ws.Range("C31").Value = 3130
$ws.Range("D31").Value = 113112
$ws.Range("E31").Value = 27341
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 44
$ws.Range("H31").Value = 5247

# Row 32: Ecuador
$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 143531
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 120511
$ws.Range("E32").Value = 11277
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 11743

# Row 42: Oman
$ws.Range("A42").Value = "Oman"
$ws.Range("B42").Value = 104129
$ws.Range("C42").Value = 664
$ws.Range("D42").Value = 91731
$ws.Range("E42").Value = 11389
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 9
$ws.Range("H42").Value = 1009

# Row 43: Egipto
$ws.Range("A43").Value = "Egipto"
$ws.Range("B43").Value = 104035
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 97492
$ws.Range("E43").Value = 533
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 6010

# Row 60: Suiza
$ws.Range("A60").Value = "Suiza"
$ws.Range("B60").Value = 58881
$ws.Range("C60").Value = 1172
$ws.Range("D60").Value = 47300
$ws.Range("E60").Value = 9497
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 2084

# Row 61: Moldavia
$ws.Range("A61").Value = "Moldavia"
$ws.Range("B61").Value = 58794
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 42480
$ws.Range("E61").Value = 14908
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1406

# Row 62: Singapur
$ws.Range("A62").Value = "Singapur"
$ws.Range("B62").Value = 57849
$ws.Range("C62").Value = 9
$ws.Range("D62").Value = 57624
$ws.Range("E62").Value = 198
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 27

# Row 92: Madagascar
$ws.Range("A92").Value = "Madagascar"
$ws.Range("B92").Value = 16654
$ws.Range("C92").Value = 21
$ws.Range("D92").Value = 15910
$ws.Range("E92").Value = 509
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 235

# Row 98: Malasia
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 14368
$ws.Range("C98").Value = 375
$ws.Range("D98").Value = 10519
$ws.Range("E98").Value = 3703
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 5
$ws.Range("H98").Value = 146

# Row 102: Finlandia
$ws.Range("A102").Value = "Finlandia"
$ws.Range("B102").Value = 11345
$ws.Range("C102").Value = 296
$ws.Range("D102").Value = 8500
$ws.Range("E102").Value = 2499
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 346

# Row 127: Hong Kong
$ws.Range("A127").Value = "Hong Kong"
$ws.Range("B127").Value = 5162
$ws.Range("C127").Value = 18
$ws.Range("D127").Value = 4890
$ws.Range("E127").Value = 167
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 105

# Row 144: Malta
$ws.Range("A144").Value = "Malta"
$ws.Range("B144").Value = 3506
$ws.Range("C144").Value = 64
$ws.Range("D144").Value = 2884
$ws.Range("E144").Value = 581
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 41

# Row 181: Gibraltar
$ws.Range("A181").Value = "Gibraltar"
$ws.Range("B181").Value = 452
$ws.Range("C181").Value = 7
$ws.Range("D181").Value = 386
$ws.Range("E181").Value = 66
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

# Row 206: Santa Lucia
$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("B206").Value = 28
$ws.Range("C206").Value = 1
$ws.Range("D206").Value = 27
$ws.Range("E206").Value = 1
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

# Row 207: Timor Oriental
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("B207").Value = 28
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 28
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

# Row 215: Montserrat
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

# Row 216: Islas Malvinas
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
